# Update countries & provincias Spain
# Applies the COVID dashboard refresh (13 Abril 2020, 06:52 -> 07:22):
#  - timestamp string bump
#  - several per-country counter updates
#  - Kirguistan's count overtakes Honduras/Reunion/Jordania/Taiwan/Malta,
#    so it re-sorts to just below Albania (rows 98-103 shift down one,
#    row 98 becomes the (updated) Kirguistan entry)
#  - Nepal's count overtakes Botsuana/Malaui/San Cristobal y Nieves,
#    so it re-sorts to just below Curazao (rows 186-189 shift down one,
#    row 186 becomes the (updated) Nepal entry)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param($row, $name, $total, $new, $active, $recovered, $critical, $deathsToday, $deaths)
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = $total
    $ws.Cells.Item($row, 3).Value = $new
    $ws.Cells.Item($row, 4).Value = $active
    $ws.Cells.Item($row, 5).Value = $recovered
    $ws.Cells.Item($row, 6).Value = $critical
    $ws.Cells.Item($row, 7).Value = $deathsToday
    $ws.Cells.Item($row, 8).Value = $deaths
}

# --- Header timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 13 de Abril de 2020 a las 07:22"

# --- Simple per-country counter refreshes (no reordering) ---
# Row 32: Australia
$ws.Cells.Item(32, 2).Value = 6351
$ws.Cells.Item(32, 3).Value = 38
$ws.Cells.Item(32, 5).Value = 2952

# Row 76: Uzbekistan
$ws.Cells.Item(76, 2).Value = 896
$ws.Cells.Item(76, 3).Value = 31
$ws.Cells.Item(76, 5).Value = 826

# Row 89: Afganistan
$ws.Cells.Item(89, 5).Value = 556
$ws.Cells.Item(89, 7).Value = 1
$ws.Cells.Item(89, 8).Value = 19

# --- Kirguistan jumps ahead of Honduras/Reunion/Jordania/Taiwan/Malta ---
# (rows 98-103 shift down by one; row 98 now holds the refreshed Kirguistan row)
Set-Row 98  "Kirguistan" 419 42 67 347 5 0 5
Set-Row 99  "Honduras"   397 4  7  365 10 0 25
Set-Row 100 "Reunion"    389 0  40 349 3  0 0
Set-Row 101 "Jordania"   389 0  201 181 5 0 7
Set-Row 102 "Taiwan"     388 0  109 273 0 0 6
Set-Row 103 "Malta"      378 0  44 331 4  0 3

# --- Nepal jumps ahead of Botsuana/Malaui/San Cristobal y Nieves ---
# (rows 186-189 shift down by one; row 186 now holds the refreshed Nepal row)
Set-Row 186 "Nepal"                   13 1 1 12 0 0 0
Set-Row 187 "Botsuana"                13 0 0 12 0 0 1
Set-Row 188 "Malaui"                  13 0 0 11 1 0 2
Set-Row 189 "San Cristobal y Nieves"  12 0 0 12 0 0 0
